# Remove the row for case 2485 (LA PLATA AV. 1095), which was merged/closed.
# All subsequent rows shift up by one, so the table shrinks from A1:P91 to A1:P90.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(34).Delete()
